# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column E corresponds to the "municipio-nombre" data field. Its curated
# metadata (rows 2-4) is updated to reflect the new dimension mapping:
#   E2: iaest-measure:municipio-nombre -> sdmx-dimension:refArea
#   E3: medida                        -> dim
#   E4: xsd:int                       -> URI-Municipio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
